# Auto-generated edit script applying numeric corrections from the
# scheduled Chocobo_Profits runner update.
$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 726.75
$ws.Range("I39").Value = 302.33334
$ws.Range("J39").Value = 2000
$ws.Range("K39").Value = 907.0000200000001
$ws.Range("L39").Value = 6000
$ws.Range("M39").Value = -611.0000200000001
$ws.Range("N39").Value = -6592
$ws.Range("H40").Value = 1500
$ws.Range("J40").Value = 1450
$ws.Range("L40").Value = 1450
$ws.Range("N40").Value = -1800
$ws.Range("H51").Value = 6830
$ws.Range("J51").Value = 7800
$ws.Range("L51").Value = 7800
$ws.Range("N51").Value = -8768
$ws.Range("H74").Value = 9094055
$ws.Range("I74").Value = 9094055
$ws.Range("K74").Value = 9094055
$ws.Range("M74").Value = -9093119
$ws.Range("H77").Value = 9094055
$ws.Range("I77").Value = 9094055
$ws.Range("K77").Value = 45470275
$ws.Range("M77").Value = -45465595
$ws.Range("H116").Value = 778424
$ws.Range("I116").Value = 1252437.5
$ws.Range("J116").Value = 20002.4
$ws.Range("K116").Value = 1252437.5
$ws.Range("L116").Value = 20002.4
$ws.Range("M116").Value = -1248995.5
$ws.Range("N116").Value = -26886.4
$ws.Range("H138").Value = 2288.15
$ws.Range("I138").Value = 1165.8387
$ws.Range("J138").Value = 2792.3767
$ws.Range("K138").Value = 3497.5161
$ws.Range("L138").Value = 8377.130099999998
$ws.Range("M138").Value = 1642.4839
$ws.Range("N138").Value = -18657.1301
$ws.Range("H141").Value = 1486.2885
$ws.Range("I141").Value = 1180.4681
$ws.Range("J141").Value = 4361
$ws.Range("K141").Value = 3541.4043
$ws.Range("L141").Value = 13083
$ws.Range("M141").Value = 1638.5957
$ws.Range("N141").Value = -23443

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9343.663
$ws.Range("I32").Value = 5941.6826
$ws.Range("J32").Value = 17586.924
$ws.Range("K32").Value = 5941.6826
$ws.Range("L32").Value = 17586.924
$ws.Range("M32").Value = -5654.6826
$ws.Range("N32").Value = -18160.924
$ws.Range("H61").Value = 1783.6842
$ws.Range("I61").Value = 1141.9166
$ws.Range("J61").Value = 2883.8572
$ws.Range("K61").Value = 1141.9166
$ws.Range("L61").Value = 2883.8572
$ws.Range("M61").Value = -929.9166
$ws.Range("N61").Value = -3307.8572
$ws.Range("H132").Value = 2039.5098
$ws.Range("I132").Value = 1065.2433
$ws.Range("K132").Value = 3195.7299
$ws.Range("M132").Value = -665.7299000000003
$ws.Range("H136").Value = 1783.6842
$ws.Range("I136").Value = 1141.9166
$ws.Range("J136").Value = 2883.8572
$ws.Range("K136").Value = 3425.7498
$ws.Range("L136").Value = 8651.571599999999
$ws.Range("M136").Value = -875.7498000000001
$ws.Range("N136").Value = -13751.5716

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2751.9207
$ws.Range("I31").Value = 1268.025
$ws.Range("J31").Value = 5332.609
$ws.Range("K31").Value = 1268.025
$ws.Range("L31").Value = 5332.609
$ws.Range("M31").Value = -973.0250000000001
$ws.Range("N31").Value = -5922.609
$ws.Range("H34").Value = 2751.9207
$ws.Range("I34").Value = 1268.025
$ws.Range("J34").Value = 5332.609
$ws.Range("K34").Value = 1268.025
$ws.Range("L34").Value = 5332.609
$ws.Range("M34").Value = -1066.025
$ws.Range("N34").Value = -5736.609

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1794.7028
$ws.Range("I5").Value = 1081.2
$ws.Range("J5").Value = 2281.182
$ws.Range("K5").Value = 3243.6
$ws.Range("L5").Value = 6843.545999999999
$ws.Range("M5").Value = -3131.6
$ws.Range("N5").Value = -7067.545999999999
$ws.Range("H17").Value = 682.0909
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H107").Value = 51057.65
$ws.Range("J107").Value = 127111
$ws.Range("L107").Value = 381333
$ws.Range("N107").Value = -385173
$ws.Range("H122").Value = 2507.0833
$ws.Range("I122").Value = 586.875
$ws.Range("J122").Value = 3467.1875
$ws.Range("K122").Value = 5281.875
$ws.Range("L122").Value = 31204.6875
$ws.Range("M122").Value = -2831.875
$ws.Range("N122").Value = -36104.6875
$ws.Range("H125").Value = 7750
$ws.Range("H126").Value = 6666.5
$ws.Range("J126").Value = 6666.5
$ws.Range("L126").Value = 19999.5
$ws.Range("N126").Value = -29879.5
$ws.Range("H130").Value = 2869.875
$ws.Range("I130").Value = 1470
$ws.Range("J130").Value = 4269.75
$ws.Range("K130").Value = 4410
$ws.Range("L130").Value = 12809.25
$ws.Range("M130").Value = 610
$ws.Range("N130").Value = -22849.25
$ws.Range("H135").Value = 1794.7028
$ws.Range("I135").Value = 1081.2
$ws.Range("J135").Value = 2281.182
$ws.Range("K135").Value = 9730.800000000001
$ws.Range("L135").Value = 20530.638
$ws.Range("M135").Value = -7195.800000000001
$ws.Range("N135").Value = -25600.638

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 30000
$ws.Range("J38").Value = 30000
$ws.Range("L38").Value = 30000
$ws.Range("N38").Value = -30926
$ws.Range("H102").Value = 2508.484
$ws.Range("I102").Value = 2095.276
$ws.Range("J102").Value = 8500
$ws.Range("K102").Value = 2095.276
$ws.Range("L102").Value = 8500
$ws.Range("M102").Value = -473.2759999999998
$ws.Range("N102").Value = -11744
$ws.Range("H127").Value = 20468.408
$ws.Range("J127").Value = 20468.408
$ws.Range("L127").Value = 20468.408
$ws.Range("N127").Value = -30388.408
$ws.Range("H130").Value = 42580
$ws.Range("J130").Value = 42580
$ws.Range("L130").Value = 42580
$ws.Range("N130").Value = -52620
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H125").Value = 37036.25
$ws.Range("I125").Value = 20000
$ws.Range("J125").Value = 42715
$ws.Range("K125").Value = 20000
$ws.Range("L125").Value = 42715
$ws.Range("M125").Value = -15080
$ws.Range("N125").Value = -52555
$ws.Range("H132").Value = 4904
$ws.Range("I132").Value = 1890.7273
$ws.Range("J132").Value = 7666.1665
$ws.Range("K132").Value = 5672.1819
$ws.Range("L132").Value = 22998.4995
$ws.Range("M132").Value = -3142.1819
$ws.Range("N132").Value = -28058.4995
$ws.Range("H134").Value = 51069.168
$ws.Range("J134").Value = 51069.168
$ws.Range("L134").Value = 51069.168
$ws.Range("N134").Value = -61209.168
$ws.Range("H136").Value = 5106.1763
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 5106.1763
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 15318.5289
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -20418.5289

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 11500
$ws.Range("I39").Value = 2600
$ws.Range("J39").Value = 15950
$ws.Range("K39").Value = 2600
$ws.Range("L39").Value = 15950
$ws.Range("M39").Value = -2187
$ws.Range("N39").Value = -16776
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H125").Value = 39730.555
$ws.Range("J125").Value = 39730.555
$ws.Range("L125").Value = 39730.555
$ws.Range("N125").Value = -49570.555
$ws.Range("H136").Value = 4533.5
$ws.Range("I136").Value = 1891.8182
$ws.Range("J136").Value = 7175.1816
$ws.Range("K136").Value = 5675.4546
$ws.Range("L136").Value = 21525.5448
$ws.Range("M136").Value = -3125.4546
$ws.Range("N136").Value = -26625.5448

